$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 182.6
$ws.Range("I17").Value = 250
$ws.Range("J17").Value = 177.13513
$ws.Range("K17").Value = 750
$ws.Range("L17").Value = 531.40539
$ws.Range("M17").Value = -582
$ws.Range("N17").Value = -867.40539

$ws.Range("H28").Value = 5716.423
$ws.Range("I28").Value = 166.89473
$ws.Range("K28").Value = 166.89473
$ws.Range("M28").Value = 318.10527

$ws.Range("H55").Value = 62500876
$ws.Range("I55").Value = 90910216
$ws.Range("J55").Value = 329.8
$ws.Range("K55").Value = 90910216
$ws.Range("L55").Value = 329.8
$ws.Range("M55").Value = -90910002
$ws.Range("N55").Value = -757.8

$ws.Range("H64").Value = 4116.6665
$ws.Range("I64").Value = 3600
$ws.Range("J64").Value = 5150
$ws.Range("K64").Value = 3600
$ws.Range("L64").Value = 5150
$ws.Range("M64").Value = -3352
$ws.Range("N64").Value = -5646

$ws.Range("H67").Value = 4116.6665
$ws.Range("I67").Value = 3600
$ws.Range("J67").Value = 5150
$ws.Range("K67").Value = 3600
$ws.Range("L67").Value = 5150
$ws.Range("M67").Value = -2742
$ws.Range("N67").Value = -6866

$ws.Range("H94").Value = 4850
$ws.Range("I94").Value = 4850
$ws.Range("K94").Value = 4850
$ws.Range("M94").Value = -4399

$ws.Range("H107").Value = 5861.364
$ws.Range("I107").Value = 6107.579
$ws.Range("J107").Value = 4302
$ws.Range("K107").Value = 6107.579
$ws.Range("L107").Value = 4302
$ws.Range("M107").Value = -4187.579
$ws.Range("N107").Value = -8142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4185.421
$ws.Range("I63").Value = 2050.0908
$ws.Range("J63").Value = 7121.5
$ws.Range("K63").Value = 2050.0908
$ws.Range("L63").Value = 7121.5
$ws.Range("M63").Value = -1364.0908
$ws.Range("N63").Value = -8493.5

$ws.Range("H66").Value = 4185.421
$ws.Range("I66").Value = 2050.0908
$ws.Range("J66").Value = 7121.5
$ws.Range("K66").Value = 10250.454
$ws.Range("L66").Value = 35607.5
$ws.Range("M66").Value = -6818.454
$ws.Range("N66").Value = -42471.5

$ws.Range("H123").Value = 27124
$ws.Range("J123").Value = 27124
$ws.Range("L123").Value = 27124
$ws.Range("N123").Value = -36924

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1850
$ws.Range("I86").Value = 1800
$ws.Range("J86").Value = 1866.6666
$ws.Range("K86").Value = 1800
$ws.Range("L86").Value = 1866.6666
$ws.Range("M86").Value = -677
$ws.Range("N86").Value = -4112.6666

$ws.Range("H89").Value = 1850
$ws.Range("I89").Value = 1800
$ws.Range("J89").Value = 1866.6666
$ws.Range("K89").Value = 9000
$ws.Range("L89").Value = 9333.333000000001
$ws.Range("M89").Value = -3384
$ws.Range("N89").Value = -20565.333

$ws.Range("H99").Value = 1016.5333
$ws.Range("I99").Value = 895.84
$ws.Range("J99").Value = 1620
$ws.Range("K99").Value = 895.84
$ws.Range("L99").Value = 1620
$ws.Range("M99").Value = 602.16
$ws.Range("N99").Value = -4616

$ws.Range("H107").Value = 1150
$ws.Range("I107").Value = 910
$ws.Range("J107").Value = 1750
$ws.Range("K107").Value = 910
$ws.Range("L107").Value = 1750
$ws.Range("M107").Value = 1010
$ws.Range("N107").Value = -5590

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 750.2
$ws.Range("I16").Value = 687.75
$ws.Range("K16").Value = 687.75
$ws.Range("M16").Value = -400.75

$ws.Range("H31").Value = 5016.951
$ws.Range("I31").Value = 1349.9333
$ws.Range("J31").Value = 7132.5386
$ws.Range("K31").Value = 1349.9333
$ws.Range("L31").Value = 7132.5386
$ws.Range("M31").Value = -1054.9333
$ws.Range("N31").Value = -7722.5386

$ws.Range("H34").Value = 5016.951
$ws.Range("I34").Value = 1349.9333
$ws.Range("J34").Value = 7132.5386
$ws.Range("K34").Value = 1349.9333
$ws.Range("L34").Value = 7132.5386
$ws.Range("M34").Value = -1147.9333
$ws.Range("N34").Value = -7536.5386

$ws.Range("H107").Value = 8929743
$ws.Range("I107").Value = 31250400
$ws.Range("J107").Value = 1480
$ws.Range("K107").Value = 31250400
$ws.Range("L107").Value = 1480
$ws.Range("M107").Value = -31248480
$ws.Range("N107").Value = -5320

$ws.Range("H109").Value = 33940
$ws.Range("J109").Value = 33940
$ws.Range("L109").Value = 33940
$ws.Range("N109").Value = -36020

$ws.Range("H113").Value = 750.2
$ws.Range("I113").Value = 687.75
$ws.Range("K113").Value = 687.75
$ws.Range("M113").Value = 1482.25

$ws.Range("H122").Value = 1953.1904
$ws.Range("I122").Value = 1803.4
$ws.Range("K122").Value = 5410.200000000001
$ws.Range("M122").Value = -2960.200000000001

$ws.Range("H134").Value = 1030.2858
$ws.Range("I134").Value = 842.4
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 2527.2
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = 7.800000000000182
$ws.Range("N134").Value = -9570

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5590.4146
$ws.Range("I70").Value = 5579.154
$ws.Range("K70").Value = 5579.154
$ws.Range("M70").Value = -5309.154

$ws.Range("H73").Value = 5590.4146
$ws.Range("I73").Value = 5579.154
$ws.Range("K73").Value = 5579.154
$ws.Range("M73").Value = -4643.154

$ws.Range("H109").Value = 10285
$ws.Range("J109").Value = 10285
$ws.Range("L109").Value = 10285
$ws.Range("N109").Value = -12365

$ws.Range("H122").Value = 5137.154
$ws.Range("I122").Value = 5334.3335
$ws.Range("K122").Value = 16003.0005
$ws.Range("M122").Value = -13553.0005

$ws.Range("H132").Value = 2207.2727
$ws.Range("I132").Value = 1609.25
$ws.Range("J132").Value = 3802
$ws.Range("K132").Value = 4827.75
$ws.Range("L132").Value = 11406
$ws.Range("M132").Value = -2297.75
$ws.Range("N132").Value = -16466

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H137").Value = 60740
$ws.Range("J137").Value = 60740
$ws.Range("L137").Value = 60740
$ws.Range("N137").Value = -70940

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws.Range("H140").Value = 60780
$ws.Range("J140").Value = 60780
$ws.Range("L140").Value = 60780
$ws.Range("N140").Value = -71140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2168.75
$ws.Range("I46").Value = 680
$ws.Range("J46").Value = 2845.4546
$ws.Range("K46").Value = 680
$ws.Range("L46").Value = 2845.4546
$ws.Range("M46").Value = -492
$ws.Range("N46").Value = -3221.4546

$ws.Range("H122").Value = 5629.2
$ws.Range("I122").Value = 3358.4
$ws.Range("J122").Value = 7900
$ws.Range("K122").Value = 10075.2
$ws.Range("L122").Value = 23700
$ws.Range("M122").Value = -7625.200000000001
$ws.Range("N122").Value = -28600

$ws.Range("H132").Value = 4248.543
$ws.Range("I132").Value = 4127.05
$ws.Range("J132").Value = 4410.533
$ws.Range("K132").Value = 12381.15
$ws.Range("L132").Value = 13231.599
$ws.Range("M132").Value = -9851.150000000001
$ws.Range("N132").Value = -18291.599

$ws.Range("H136").Value = 9805683
$ws.Range("I136").Value = 2152
$ws.Range("J136").Value = 18519934
$ws.Range("K136").Value = 6456
$ws.Range("L136").Value = 55559802
$ws.Range("M136").Value = -3906
$ws.Range("N136").Value = -55564902

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5474.3335
$ws.Range("I81").Value = 5711.6665
$ws.Range("J81").Value = 4999.6665
$ws.Range("K81").Value = 11423.333
$ws.Range("L81").Value = 9999.333000000001
$ws.Range("M81").Value = -10362.333
$ws.Range("N81").Value = -12121.333

$ws.Range("H84").Value = 5474.3335
$ws.Range("I84").Value = 5711.6665
$ws.Range("J84").Value = 4999.6665
$ws.Range("K84").Value = 57116.665
$ws.Range("L84").Value = 49996.665
$ws.Range("M84").Value = -51812.665
$ws.Range("N84").Value = -60604.665

$ws.Range("H96").Value = 6356
$ws.Range("J96").Value = 7866.6665
$ws.Range("L96").Value = 7866.6665
$ws.Range("N96").Value = -10612.6665

$ws.Range("H100").Value = 850
$ws.Range("I100").Value = 673.3333
$ws.Range("J100").Value = 1026.6666
$ws.Range("K100").Value = 1346.6666
$ws.Range("L100").Value = 2053.3332
$ws.Range("M100").Value = -805.6666
$ws.Range("N100").Value = -3135.3332

$ws.Range("H123").Value = 24672.357
$ws.Range("J123").Value = 24672.357
$ws.Range("L123").Value = 24672.357
$ws.Range("N123").Value = -34472.357
